$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.922.38"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").Value = "2.997.14"
$ws.Range("E3").Value = "  +2.30%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.28"
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.63"
$ws.Range("E6").Value = "  -3.92%  "

$ws.Range("E7").Value = "  -2.14%  "

$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").Value = "  -3.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.94"
$ws.Range("E10").Value = "  -3.87%  "

$ws.Range("E11").Value = "  +2.68%  "

$ws.Range("E12").Value = "  -3.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.00"
$ws.Range("E13").Value = "  -3.69%  "

$ws.Range("D14").Value = "3.476.46"
$ws.Range("E14").Value = "  +2.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.57"
$ws.Range("E15").Value = "  -4.15%  "

$ws.Range("D16").Value = "2.983.75"
$ws.Range("E16").Value = "  +2.24%  "

$ws.Range("E17").Value = "  +1.44%  "

$ws.Range("D18").Value = "51.910.45"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.38"
$ws.Range("E19").Value = "  +2.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.42"
$ws.Range("E20").Value = "  -2.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.46"
$ws.Range("E21").Value = "  -3.97%  "

$ws.Range("D22").Value = "0.0₃0967"
$ws.Range("E22").Value = "  -1.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.01"
$ws.Range("E23").Value = "  -2.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.22"
$ws.Range("E24").Value = "  -2.85%  "

$ws.Range("E25").Value = "  -4.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.178"
$ws.Range("E26").Value = "  -3.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.84"
$ws.Range("E27").Value = "  -1.35%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.35"
$ws.Range("E29").Value = "  -0.64%  "

$ws.Range("E30").Value = "  +2.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.31"
$ws.Range("E31").Value = "  +4.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.13"
$ws.Range("E32").Value = "  -4.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.99"
$ws.Range("E33").Value = "  -7.64%  "

$ws.Range("E34").Value = "  +12.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.03"
$ws.Range("E35").Value = "  -2.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0430"
$ws.Range("E36").Value = "  -3.30%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.26"
$ws.Range("E38").Value = "  +0.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.83"
$ws.Range("E39").Value = "  +3.08%  "

$ws.Range("E40").Value = "  -3.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.44"
$ws.Range("E41").Value = "  -6.00%  "

$ws.Range("E42").Value = "  -3.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.14"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("E44").Value = "  +3.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.16"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("D46").Value = "2.119.07"
$ws.Range("E46").Value = "  -1.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.30"
$ws.Range("E47").Value = "  -4.76%  "

$ws.Range("E48").Value = "  -7.56%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "3.300.27"
$ws.Range("E49").Value = "  +2.56%  "

$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.242"
$ws.Range("E50").Value = "  -2.94%  "

$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0331"
